$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.261679768562317
$ws.Range("B1").Value = 2.331920862197876
$ws.Range("C1").Value = 4.494894027709961
$ws.Range("D1").Value = 2.843429565429688
$ws.Range("E1").Value = 1.355108499526978
